$wb = $excel.ActiveWorkbook

# --- Data sheet (row 2 / row 3): stock test case updated for 0 quantity ---
$ws = $wb.Worksheets.Item("Data")

# Row 2
$ws.Cells.Item(2, 1).Value = "Transfer out Imprest"      # A2: Transfer in Imprest -> Transfer out Imprest
$ws.Cells.Item(2, 3).Value = "Acetec 5 mg tablet"         # C2: quinapril 10 mg tablet -> Acetec 5 mg tablet
$ws.Cells.Item(2, 5).Value = "Perry Grant"                # E2: Anthony Jones -> Perry Grant

# Row 3
$ws.Cells.Item(3, 3).Value = "Endone 5 mg tablet"         # C3: (prilocaine) ... -> Endone 5 mg tablet
$ws.Cells.Item(3, 4).Clear()                              # D3: drop the empty/styled cell entirely
$ws.Cells.Item(3, 5).Value = "Perry Grant"                # E3: Matilda Kerr -> Perry Grant

# --- Cursor/selection bookkeeping left behind by the edit session ---
$cfg = $wb.Worksheets.Item("Configuration")
$cfg.Range("D2").Select()

$ws.Range("C5").Select()
